$wb = $excel.ActiveWorkbook

# F2:F9 values to apply (想去人数 / "want to go" counts)
$values = @(354, 96, 1538, 22, 51, 134, 56, 364)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 6).Value = $values[$i]
    }
}
